$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the StudentTerm sheet: append the new term rows (19-28) for
#    student S533622, rows 20-29 of the sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("StudentTerm")

$rows = @(
    @(19, "S533622", 1, "F18",  "FALL2018"),
    @(20, "S533622", 2, "S19",  "Spring2019"),
    @(21, "S533622", 3, "Sum19","Summer2019"),
    @(22, "S533622", 4, "F19",  "Fall2019"),
    @(23, "S533622", 5, "S20",  "Spring2020"),
    @(24, "S533622", 1, "F18",  "FALL2018"),
    @(25, "S533622", 2, "S19",  "Spring2019"),
    @(26, "S533622", 3, "F19",  "Fall2019"),
    @(27, "S533622", 4, "S20",  "Spring2020"),
    @(28, "S533622", 5, "Sum20","Summer2020")
)

$r = 20
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Row 19 (the former last row of the block) loses its bottom border so it
# visually merges with the new block of rows that follows it.
$ws.Range("A19:E19").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

# Rows 20-29 keep the boxed border, but switch to a plain (non-bold) black
# font and wrap their text.
$dataRange = $ws.Range("A20:E29")
$dataRange.Font.Bold = $false
$dataRange.Font.Color = 0
$dataRange.WrapText = $true

$ws.Range("A20:A29").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C20:C29").HorizontalAlignment = -4152  # xlRight

# ---------------------------------------------------------------------------
# 2. Switch the active sheet/selection: StudentTerm becomes the active tab
#    (it was DegreePlan before), with a new selection on StudentTerm.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("L22").Select() | Out-Null
